$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13 and 14 need to have their contents swapped (columns A through AY).
# Column Y (25) and AA (27) hold dates stored as plain text (e.g. "2023-08-12").
# Excel auto-converts such text to a real date serial number when assigned
# through .Value, so those two columns must be forced back to text format
# before the value is written, to preserve the original text representation.
$lastCol = 51  # column AY
$textColumns = @(25, 27)

for ($col = 1; $col -le $lastCol; $col++) {
    $cell13 = $ws.Cells.Item(13, $col)
    $cell14 = $ws.Cells.Item(14, $col)

    $value13 = $cell13.Value()
    $value14 = $cell14.Value()

    # Skip columns that are blank on both rows so we don't disturb
    # already-empty placeholder cells.
    $empty13 = ($value13 -eq $null) -or ($value13 -eq "")
    $empty14 = ($value14 -eq $null) -or ($value14 -eq "")
    if ($empty13 -and $empty14) {
        continue
    }

    if ($textColumns -contains $col) {
        $cell13.NumberFormat = "@"
        $cell14.NumberFormat = "@"
    }

    $cell13.Value = $value14
    $cell14.Value = $value13
}
